$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D (shifts existing D:K data to E:L,
# and the empty K column to L) to make room for the new fiscal-year (FY2018)
# figures that were added as the first data column of each table.
$ws.Columns("D").Insert()

# The freshly inserted column D has no number formatting yet; copy the
# formatting (date format / number format) from column E - which now
# holds what used to be column D - across the three tables so the new
# column matches its neighbours exactly.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)

$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)

$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Populate the new column D with the FY2018 figures for the Income
# Statement, Balance Sheet and Cash Flow Statement tables.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 1031600
$ws.Range("D9").Value = 367700
$ws.Range("D10").Value = 663900
$ws.Range("D12").Value = 417400
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1217900
$ws.Range("D18").Value = -186300
$ws.Range("D20").Value = 82500
$ws.Range("D21").Value = -74700
$ws.Range("D22").Value = 22700
$ws.Range("D23").Value = -126500
$ws.Range("D24").Value = 1000
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -127500
$ws.Range("D27").Value = -127500
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 400
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -82500
$ws.Range("D33").Value = -127100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -127100
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 1137000
$ws.Range("D42").Value = 248600
$ws.Range("D43").Value = 226700
$ws.Range("D44").Value = 70700
$ws.Range("D45").Value = 16500
$ws.Range("D46").Value = 1699500
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 183100
$ws.Range("D49").Value = 18700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 14700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 1916000
$ws.Range("D57").Value = 75500
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 146900
$ws.Range("D60").Value = 222400
$ws.Range("D61").Value = 1017600
$ws.Range("D62").Value = 12700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1252700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -798900
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 663300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = -127100
$ws.Range("D83").Value = 29100
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 123200
$ws.Range("D91").Value = -67100
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -139800
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 710400
$ws.Range("D101").Value = 1800
$ws.Range("D102").Value = 695600
